$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1070.5
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 144
$ws.Cells.Item(2, 6).Value = 3652.23
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 4866.73

# Row 3
$ws.Cells.Item(3, 2).Value = 1129.46
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 288
$ws.Cells.Item(3, 6).Value = 3693.63
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 5111.09

# Row 4
$ws.Cells.Item(4, 2).Value = 1165.75
$ws.Cells.Item(4, 3).Value = 82.64
$ws.Cells.Item(4, 4).Value = 128.86
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3581.92
$ws.Cells.Item(4, 7).Value = 3.39
$ws.Cells.Item(4, 8).Value = 4962.56

# Row 5
$ws.Cells.Item(5, 2).Value = 1687.39
$ws.Cells.Item(5, 3).Value = 1198.18
$ws.Cells.Item(5, 4).Value = 1443.77
$ws.Cells.Item(5, 5).Value = 432
$ws.Cells.Item(5, 6).Value = 4437.14
$ws.Cells.Item(5, 7).Value = 289.47
$ws.Cells.Item(5, 8).Value = 9487.949999999999

# Row 6
$ws.Cells.Item(6, 2).Value = 4331.88
$ws.Cells.Item(6, 3).Value = 5527.2
$ws.Cells.Item(6, 4).Value = 6083.98
$ws.Cells.Item(6, 5).Value = 432
$ws.Cells.Item(6, 6).Value = 10224.64
$ws.Cells.Item(6, 7).Value = 2599.44
$ws.Cells.Item(6, 8).Value = 29199.139999999996

# Row 7
$ws.Cells.Item(7, 2).Value = 8164.8
$ws.Cells.Item(7, 3).Value = 11322.69
$ws.Cells.Item(7, 4).Value = 12906.9
$ws.Cells.Item(7, 5).Value = 2160
$ws.Cells.Item(7, 6).Value = 29066.32
$ws.Cells.Item(7, 7).Value = 11761.73
$ws.Cells.Item(7, 8).Value = 75382.44

# Row 8
$ws.Cells.Item(8, 2).Value = 14633.14
$ws.Cells.Item(8, 3).Value = 25971.77
$ws.Cells.Item(8, 4).Value = 24584.07
$ws.Cells.Item(8, 5).Value = 5472
$ws.Cells.Item(8, 6).Value = 59102.82
$ws.Cells.Item(8, 7).Value = 28363.75
$ws.Cells.Item(8, 8).Value = 158127.55000000002

# Row 9
$ws.Cells.Item(9, 2).Value = 16270.63
$ws.Cells.Item(9, 3).Value = 35951.54
$ws.Cells.Item(9, 4).Value = 30924
$ws.Cells.Item(9, 5).Value = 5328
$ws.Cells.Item(9, 6).Value = 41854.59
$ws.Cells.Item(9, 7).Value = 17046.15
$ws.Cells.Item(9, 8).Value = 147374.91

# Row 10
$ws.Cells.Item(10, 2).Value = 4309.2
$ws.Cells.Item(10, 3).Value = 8375.99
$ws.Cells.Item(10, 4).Value = 7987.75
$ws.Cells.Item(10, 5).Value = 1152
$ws.Cells.Item(10, 6).Value = 11374.67
$ws.Cells.Item(10, 7).Value = 3667.27
$ws.Cells.Item(10, 8).Value = 36866.88

# Row 11
$ws.Cells.Item(11, 2).Value = 2159.14
$ws.Cells.Item(11, 3).Value = 3700.21
$ws.Cells.Item(11, 4).Value = 3570.17
$ws.Cells.Item(11, 5).Value = 864
$ws.Cells.Item(11, 6).Value = 6049.94
$ws.Cells.Item(11, 7).Value = 2548.78
$ws.Cells.Item(11, 8).Value = 18892.239999999998

# Row 12
$ws.Cells.Item(12, 2).Value = 1628.42
$ws.Cells.Item(12, 3).Value = 2614.04
$ws.Cells.Item(12, 4).Value = 2526.46
$ws.Cells.Item(12, 5).Value = 144
$ws.Cells.Item(12, 6).Value = 4622.61
$ws.Cells.Item(12, 7).Value = 1480.87
$ws.Cells.Item(12, 8).Value = 13016.399999999998

# Row 13
$ws.Cells.Item(13, 2).Value = 2503.87
$ws.Cells.Item(13, 3).Value = 4118.1
$ws.Cells.Item(13, 4).Value = 4727.63
$ws.Cells.Item(13, 5).Value = 720
$ws.Cells.Item(13, 6).Value = 4122.06
$ws.Cells.Item(13, 7).Value = 786.04
$ws.Cells.Item(13, 8).Value = 16977.7

# Row 14
$ws.Cells.Item(14, 2).Value = 4272.91
$ws.Cells.Item(14, 3).Value = 6657.7
$ws.Cells.Item(14, 4).Value = 11189.88
$ws.Cells.Item(14, 5).Value = 432
$ws.Cells.Item(14, 6).Value = 5179.66
$ws.Cells.Item(14, 7).Value = 552.49
$ws.Cells.Item(14, 8).Value = 28284.64

# Row 15
$ws.Cells.Item(15, 2).Value = 6391.22
$ws.Cells.Item(15, 3).Value = 5767.98
$ws.Cells.Item(15, 4).Value = 15188.5
$ws.Cells.Item(15, 5).Value = 864
$ws.Cells.Item(15, 6).Value = 9357.3
$ws.Cells.Item(15, 7).Value = 185.13
$ws.Cells.Item(15, 8).Value = 37754.13

# Row 16
$ws.Cells.Item(16, 2).Value = 7702.13
$ws.Cells.Item(16, 3).Value = 6155.52
$ws.Cells.Item(16, 4).Value = 15685.31
$ws.Cells.Item(16, 5).Value = 1584
$ws.Cells.Item(16, 6).Value = 17784.25
$ws.Cells.Item(16, 7).Value = 150.1
$ws.Cells.Item(16, 8).Value = 49061.31

# Row 17
$ws.Cells.Item(17, 2).Value = 10500.84
$ws.Cells.Item(17, 3).Value = 9217.81
$ws.Cells.Item(17, 4).Value = 16667.26
$ws.Cells.Item(17, 5).Value = 1728
$ws.Cells.Item(17, 6).Value = 30091.36
$ws.Cells.Item(17, 7).Value = 34.39
$ws.Cells.Item(17, 8).Value = 68239.66

# Row 18
$ws.Cells.Item(18, 2).Value = 12836.88
$ws.Cells.Item(18, 3).Value = 11381.7
$ws.Cells.Item(18, 4).Value = 14356.91
$ws.Cells.Item(18, 5).Value = 3312
$ws.Cells.Item(18, 6).Value = 45621
$ws.Cells.Item(18, 7).Value = 13.42
$ws.Cells.Item(18, 8).Value = 87521.91

# Row 19
$ws.Cells.Item(19, 2).Value = 14909.83
$ws.Cells.Item(19, 3).Value = 14494.58
$ws.Cells.Item(19, 4).Value = 11755.61
$ws.Cells.Item(19, 5).Value = 4320
$ws.Cells.Item(19, 6).Value = 60054.78
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 105534.8

# Row 20
$ws.Cells.Item(20, 2).Value = 16697.02
$ws.Cells.Item(20, 3).Value = 23225.98
$ws.Cells.Item(20, 4).Value = 14902.16
$ws.Cells.Item(20, 5).Value = 4320
$ws.Cells.Item(20, 6).Value = 67508.29
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 126653.45

# Row 21
$ws.Cells.Item(21, 2).Value = 16361.35
$ws.Cells.Item(21, 3).Value = 18145.71
$ws.Cells.Item(21, 4).Value = 12150.66
$ws.Cells.Item(21, 5).Value = 4320
$ws.Cells.Item(21, 6).Value = 65763.73
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 116741.45

# Row 22
$ws.Cells.Item(22, 2).Value = 13653.36
$ws.Cells.Item(22, 3).Value = 13465.59
$ws.Cells.Item(22, 4).Value = 11924.95
$ws.Cells.Item(22, 5).Value = 5040
$ws.Cells.Item(22, 6).Value = 56300.88
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 100384.78

# Row 23
$ws.Cells.Item(23, 2).Value = 11984.11
$ws.Cells.Item(23, 3).Value = 9547.84
$ws.Cells.Item(23, 4).Value = 14809.43
$ws.Cells.Item(23, 5).Value = 2880
$ws.Cells.Item(23, 6).Value = 43801.83
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 83023.21

# Row 24
$ws.Cells.Item(24, 2).Value = 8359.85
$ws.Cells.Item(24, 3).Value = 6859.34
$ws.Cells.Item(24, 4).Value = 12564.61
$ws.Cells.Item(24, 5).Value = 2448
$ws.Cells.Item(24, 6).Value = 26128.34
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 56360.14

# Row 25
$ws.Cells.Item(25, 2).Value = 5470.42
$ws.Cells.Item(25, 3).Value = 4105.49
$ws.Cells.Item(25, 4).Value = 9834.86
$ws.Cells.Item(25, 5).Value = 1296
$ws.Cells.Item(25, 6).Value = 13119.88
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 33826.65

# Row 26
$ws.Cells.Item(26, 2).Value = 188194.10000000003
$ws.Cells.Item(26, 3).Value = 227887.59999999998
$ws.Cells.Item(26, 4).Value = 255913.73000000004
$ws.Cells.Item(26, 5).Value = 49680
$ws.Cells.Item(26, 6).Value = 622493.8699999999
$ws.Cells.Item(26, 7).Value = 69482.42
$ws.Cells.Item(26, 8).Value = 1413651.7199999997

